$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.617.46"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "2.091.97"
$ws.Range("E3").Value = "  +9.32%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.73"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.660"
$ws.Range("E6").Value = "  -5.82%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.37"
$ws.Range("E8").Value = "  +5.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.35"
$ws.Range("E9").Value = "  +3.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.376"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0747"
$ws.Range("E11").Value = "  -1.75%  "
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.86"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("D14").Value = "2.399.08"
$ws.Range("E14").Value = "  +9.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.839"
$ws.Range("E15").Value = "  +2.70%  "
$ws.Range("D16").Value = "2.095.14"
$ws.Range("E16").Value = "  +9.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.13"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").Value = "36.581.30"
$ws.Range("E18").Value = "  -1.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.23"
$ws.Range("E19").Value = "  -1.96%  "
$ws.Range("D20").Value = "0.0₃0837"
$ws.Range("E20").Value = "  -2.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.30"
$ws.Range("E21").Value = "  -2.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "240.50"
$ws.Range("E22").Value = "  -4.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.33"
$ws.Range("E23").Value = "  +2.91%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.54"
$ws.Range("E25").Value = "  -2.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.39"
$ws.Range("E26").Value = "  +2.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.29"
$ws.Range("E27").Value = "  +5.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.18"
$ws.Range("E28").Value = "  +13.09%  "
$ws.Range("E29").Value = "  -10.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.47"
$ws.Range("E30").Value = "  +45.65%  "
$ws.Range("E31").Value = "  -5.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.51"
$ws.Range("E32").Value = "  -2.08%  "
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +13.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.43"
$ws.Range("E35").Value = "  +20.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0910"
$ws.Range("E36").Value = "  +2.63%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  -3.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.10"
$ws.Range("E39").Value = "  -5.43%  "
$ws.Range("E40").Value = "  -9.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0225"
$ws.Range("E41").Value = "  -1.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.17"
$ws.Range("E42").Value = "  +6.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "98.13"
$ws.Range("E43").Value = "  -6.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.64"
$ws.Range("E44").Value = "  -7.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.78"
$ws.Range("E45").Value = "  -3.90%  "
$ws.Range("D46").Value = "1.337.73"
$ws.Range("E46").Value = "  -0.99%  "
$ws.Range("E47").Value = "  +3.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.07"
$ws.Range("E48").Value = "  +9.05%  "
$ws.Range("D50").Value = "2.285.00"
$ws.Range("E50").Value = "  +9.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.27"
$ws.Range("E51").Value = "  -5.94%  "
